# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates the DAMSLTag (I) and DialogAct (J)
# columns for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 12;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 13;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 15;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 23;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 24;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 26;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 28;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 34;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 76;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 102; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 105; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 109; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 113; Tag = "sv"; Act = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
